# CS133JS_Lab06_Rubric.xlsx - "New and updated notes"
#
# Restructure the "Part 2" rubric rows (10-14) so their criteria text sits
# in column B instead of column A (freeing column A as a narrow left
# margin/indent column), and size the columns to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rubric")

# --- Move the Part 2 criteria labels from column A to column B ------------
# Rows 10-14 currently hold their text in A10:A14; shift it one column right.
$ws.Range("A10:A14").Cut($ws.Range("B10:B14")) | Out-Null

# --- Column widths for the new layout --------------------------------------
# Column A becomes a thin indent/margin column.
$ws.Columns.Item(1).ColumnWidth = 1.8333333333333333
# Column C (criteria text) widens to fit the longer wrapped labels.
$ws.Columns.Item(3).ColumnWidth = 36
# Columns D/E (Possible / Score) get slightly narrower, fitted widths.
$ws.Columns.Item(4).ColumnWidth = 6.833333333333333
$ws.Columns.Item(5).ColumnWidth = 5
# Column F becomes a thin trailing margin column.
$ws.Columns.Item(6).ColumnWidth = 0.6666666666666666

# --- Update the active selection -------------------------------------------
$ws.Activate()
$ws.Range("G13").Select() | Out-Null
